# Update "PM12 Tidsregistrering for Emil.xlsx" - add a new time-registration
# entry (row 27) on the "Ark1" sheet, then move the view/selection as the
# author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")
$ws.Activate()

# --- Fill in the new row 27 entry -----------------------------------------
# A27: task description (new shared string)
$ws.Range("A27").Value = "Lav FXML UI til UC06 kontante kapacitetsomkostninger"
# B27: role (reuses existing "User-Interface Designer" shared string)
$ws.Range("B27").Value = "User-Interface Designer"
# C27: date (06-03-2020)
$ws.Range("C27").Value = 43896
# D27: start time (10:00)
$ws.Range("D27").Value = 0.41666666666666669
# E27: end time (13:10)
$ws.Range("E27").Value = 0.54861111111111105
# F27: duration (02:00) - copy formatting from the row above since F27 had
# no explicit formatting yet (it only inherited the bare column style)
$ws.Range("F27").Value = 0.083333333333333329
$ws.Range("F26").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# G27 and H27 already contain formulas (shared formula / SUM) that will
# recalculate automatically once D27:F27 have values.
$excel.Calculate()

# --- Update view state ------------------------------------------------------
$ws.Range("F33").Select()
